# Update the "想去人数" (people interested) counts in the "展览" and
# "全部类型" worksheets, per the source commit's regenerated data output.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1089
$ws1.Range("F5").Value  = 425
$ws1.Range("F7").Value  = 551
$ws1.Range("F8").Value  = 65
$ws1.Range("F9").Value  = 6775
$ws1.Range("F16").Value = 16155
$ws1.Range("F17").Value = 1584
$ws1.Range("F19").Value = 328
$ws1.Range("F22").Value = 11329
$ws1.Range("F24").Value = 964
$ws1.Range("F25").Value = 4455
$ws1.Range("F29").Value = 42
$ws1.Range("F30").Value = 318
$ws1.Range("F31").Value = 139

# ---- Sheet "全部类型" ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 1089
$ws4.Range("F5").Value  = 425
$ws4.Range("F7").Value  = 551
$ws4.Range("F9").Value  = 65
$ws4.Range("F10").Value = 6775
$ws4.Range("F18").Value = 16155
$ws4.Range("F19").Value = 1584
$ws4.Range("F21").Value = 328
$ws4.Range("F26").Value = 11329
$ws4.Range("F28").Value = 964
$ws4.Range("F29").Value = 4455
$ws4.Range("F33").Value = 42
$ws4.Range("F34").Value = 318
$ws4.Range("F35").Value = 139

$wb.Save()
